$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33500
$ws.Range("J3").Value = 33500
$ws.Range("L3").Value = 33500
$ws.Range("N3").Value = -33728

$ws.Range("H5").Value = 147.09091
$ws.Range("I5").Value = 119.77778
$ws.Range("K5").Value = 119.77778
$ws.Range("M5").Value = -4.777780000000007

$ws.Range("H7").Value = 8775
$ws.Range("J7").Value = 8775
$ws.Range("L7").Value = 8775
$ws.Range("N7").Value = -8999

$ws.Range("H14").Value = 8775
$ws.Range("J14").Value = 8775
$ws.Range("L14").Value = 8775
$ws.Range("N14").Value = -9157

$ws.Range("H102").Value = 33500
$ws.Range("J102").Value = 33500
$ws.Range("L102").Value = 33500
$ws.Range("N102").Value = -39990

$ws.Range("H106").Value = 4986
$ws.Range("I106").Value = 4986
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4986
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -4355
$ws.Range("N106").ClearContents()

$ws.Range("H135").Value = 1662.05
$ws.Range("I135").Value = 1736.421
$ws.Range("K135").Value = 15627.789
$ws.Range("M135").Value = -13092.789

$ws.Range("H137").Value = 12991.15
$ws.Range("I137").Value = 1832.0769
$ws.Range("K137").Value = 5496.2307
$ws.Range("M137").Value = -2946.2307

$ws.Range("H138").Value = 2136.7976
$ws.Range("J138").Value = 2361.0344
$ws.Range("L138").Value = 7083.1032
$ws.Range("N138").Value = -17363.1032

$ws.Range("H141").Value = 3585.375
$ws.Range("I141").Value = 3950.8462
$ws.Range("K141").Value = 11852.5386
$ws.Range("M141").Value = -6672.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 144.57143
$ws.Range("J4").Value = 995
$ws.Range("L4").Value = 995
$ws.Range("N4").Value = -1227

$ws.Range("H43").Value = 31499.5
$ws.Range("I43").Value = 29999
$ws.Range("K43").Value = 29999
$ws.Range("M43").Value = -29686

$ws.Range("H61").Value = 55927.953
$ws.Range("I61").Value = 1357.3103
$ws.Range("K61").Value = 1357.3103
$ws.Range("M61").Value = -1145.3103

$ws.Range("H102").Value = 6494.143
$ws.Range("I102").Value = 4491.8
$ws.Range("J102").Value = 11500
$ws.Range("K102").Value = 4491.8
$ws.Range("L102").Value = 11500
$ws.Range("M102").Value = -2869.8
$ws.Range("N102").Value = -14744

$ws.Range("H136").Value = 55927.953
$ws.Range("I136").Value = 1357.3103
$ws.Range("K136").Value = 4071.9309
$ws.Range("M136").Value = -1521.9309

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 20690.959
$ws.Range("I99").Value = 24023.9
$ws.Range("J99").Value = 4026.25
$ws.Range("K99").Value = 24023.9
$ws.Range("L99").Value = 4026.25
$ws.Range("M99").Value = -22525.9
$ws.Range("N99").Value = -7022.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 975.12
$ws.Range("I7").Value = 802.26666
$ws.Range("J7").Value = 1234.4
$ws.Range("K7").Value = 802.26666
$ws.Range("L7").Value = 1234.4
$ws.Range("M7").Value = -689.26666
$ws.Range("N7").Value = -1460.4

$ws.Range("H54").Value = 14296
$ws.Range("I54").Value = 14494
$ws.Range("J54").Value = 14246.5
$ws.Range("K54").Value = 14494
$ws.Range("L54").Value = 14246.5
$ws.Range("M54").Value = -13836
$ws.Range("N54").Value = -15562.5

$ws.Range("H58").Value = 10700.125
$ws.Range("I58").Value = 5117.6333
$ws.Range("K58").Value = 5117.6333
$ws.Range("M58").Value = -4914.6333

$ws.Range("H88").Value = 58210.5
$ws.Range("J88").Value = 58210.5
$ws.Range("L88").Value = 58210.5
$ws.Range("N88").Value = -59022.5

$ws.Range("H91").Value = 58210.5
$ws.Range("J91").Value = 58210.5
$ws.Range("L91").Value = 58210.5
$ws.Range("N91").Value = -61018.5

$ws.Range("H134").Value = 26320920
$ws.Range("I134").Value = 1286.862
$ws.Range("J134").Value = 111128620
$ws.Range("K134").Value = 3860.586
$ws.Range("L134").Value = 333385860
$ws.Range("M134").Value = -1325.586
$ws.Range("N134").Value = -333390930

$ws.Range("H136").Value = 10700.125
$ws.Range("I136").Value = 5117.6333
$ws.Range("K136").Value = 15352.8999
$ws.Range("M136").Value = -12802.8999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.5
$ws.Range("I2").Value = 147.14285
$ws.Range("K2").Value = 882.8571000000001
$ws.Range("M2").Value = -769.8571000000001

$ws.Range("H4").Value = 2135342.2
$ws.Range("I4").Value = 2862053.5
$ws.Range("K4").Value = 8586160.5
$ws.Range("M4").Value = -8586048.5

$ws.Range("H24").Value = 1599
$ws.Range("J24").Value = 1599
$ws.Range("L24").Value = 4797
$ws.Range("N24").Value = -5257

$ws.Range("H68").Value = 4347.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4347.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 13042.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -14664.5

$ws.Range("H71").Value = 4347.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4347.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 39127.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -47239.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3092.5454
$ws.Range("I61").Value = 2377.375
$ws.Range("J61").Value = 4999.6665
$ws.Range("K61").Value = 2377.375
$ws.Range("L61").Value = 4999.6665
$ws.Range("M61").Value = -2175.375
$ws.Range("N61").Value = -5403.6665

$ws.Range("H93").Value = 50005290
$ws.Range("I93").Value = 125007080
$ws.Range("K93").Value = 125007080
$ws.Range("M93").Value = -125005832

$ws.Range("H100").Value = 3794.3
$ws.Range("I100").Value = 3991.8572
$ws.Range("J100").Value = 3333.3333
$ws.Range("K100").Value = 3991.8572
$ws.Range("L100").Value = 3333.3333
$ws.Range("M100").Value = -3450.8572
$ws.Range("N100").Value = -4415.3333

$ws.Range("H113").Value = 3092.5454
$ws.Range("I113").Value = 2377.375
$ws.Range("J113").Value = 4999.6665
$ws.Range("K113").Value = 2377.375
$ws.Range("L113").Value = 4999.6665
$ws.Range("M113").Value = -207.375
$ws.Range("N113").Value = -9339.666499999999

$ws.Range("H122").Value = 27291372
$ws.Range("I122").Value = 41662376
$ws.Range("J122").Value = 3339698.8
$ws.Range("K122").Value = 124987128
$ws.Range("L122").Value = 10019096.4
$ws.Range("M122").Value = -124984678
$ws.Range("N122").Value = -10023996.4

$ws.Range("H132").Value = 1833889.9
$ws.Range("I132").Value = 9407.5
$ws.Range("K132").Value = 28222.5
$ws.Range("M132").Value = -25692.5

$ws.Range("H136").Value = 9522.924999999999
$ws.Range("I136").Value = 6962.1797
$ws.Range("K136").Value = 20886.5391
$ws.Range("M136").Value = -18336.5391

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 643.8570999999999
$ws.Range("I100").Value = 562.7143
$ws.Range("J100").Value = 725
$ws.Range("K100").Value = 1125.4286
$ws.Range("L100").Value = 1450
$ws.Range("M100").Value = -584.4286
$ws.Range("N100").Value = -2532

$ws.Range("H107").Value = 1135.375
$ws.Range("I107").Value = 1182
$ws.Range("K107").Value = 3546
$ws.Range("M107").Value = -1626

$ws.Range("H126").Value = 6254415
$ws.Range("I126").Value = 5261.3076
$ws.Range("J126").Value = 33334080
$ws.Range("K126").Value = 15783.9228
$ws.Range("L126").Value = 100002240
$ws.Range("M126").Value = -13313.9228
$ws.Range("N126").Value = -100007180

$ws.Range("H132").Value = 27752.428
$ws.Range("J132").Value = 67817.336
$ws.Range("L132").Value = 203452.008
$ws.Range("N132").Value = -208512.008

$ws.Range("H135").Value = 100715
$ws.Range("J135").Value = 100715
$ws.Range("L135").Value = 100715
$ws.Range("N135").Value = -110855

$ws.Range("H136").Value = 14474.412
$ws.Range("I136").Value = 2293.2
$ws.Range("K136").Value = 6879.599999999999
$ws.Range("M136").Value = -4329.599999999999
